$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bearing choices")

# Update bearing width dimension (column L) with the latest measurements
# from CAD. Rows 3-7 share one bearing width; row 8 uses a different
# bearing. Column M ("SF = C0 / calc load") is a shared formula (K/L) and
# recalculates automatically.
$ws.Range("L3:L7").Value = 9.9644999999999992
$ws.Range("L8").Value = 2.3117999999999999

# Bring the updated columns into view and leave the selection where the
# author left it when saving.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 6
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("L6").Select()
